# The template has a field whose field code is " m:1/0 " (a divide-by-zero
# M2Doc query). The commit turns the field construct (fldChar begin /
# instrText.../ fldChar end runs) into plain literal text runs that just
# spell out the same token stream as "{", "m", ":1/0", "}" - i.e. the field
# delimiters/instrText are rewritten as ordinary <w:t> runs instead of a
# live Word field.

$d = $word.ActiveDocument

# Locate the field and the paragraph that hosts it (field code / result
# ranges collapse when the result is empty, so we resolve the owning
# paragraph from the field's Result position rather than assuming a fixed
# paragraph index).
$f = $d.Fields.Item(1)
$anchor = $f.Result.Start
$paraIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($anchor -ge $candidate.Range.Start -and $anchor -le $candidate.Range.End) {
        $paraIndex = $i
        break
    }
}

# Remove the field (fldChar begin/end + instrText runs) entirely, leaving
# an empty paragraph behind in its place.
$f.Delete()

# Re-type the field's token stream as plain text, as four separate runs:
# "{", "m", ":1/0", "}". Each chunk is inserted through a transient
# bookmark at the insertion point so the new run does not get silently
# coalesced into the previous one (same formatting would otherwise merge
# them into a single run); the bookmark itself is deleted right after so
# no bookmark markup is left behind in the document.
$chunks = @("{", "m", ":1/0", "}")
$markIndex = 0
foreach ($chunk in $chunks) {
    $p = $d.Paragraphs.Item($paraIndex)
    $insertAt = $d.Range($p.Range.End - 1, $p.Range.End - 1)

    if ($insertAt.Start -gt $p.Range.Start) {
        $markIndex = $markIndex + 1
        $markName = "m2doc_fieldrewrite_$markIndex"
        $d.Bookmarks.Add($markName, $insertAt) | Out-Null
        $insertAt.InsertAfter($chunk)
        $d.Bookmarks($markName).Delete()
    } else {
        $insertAt.InsertAfter($chunk)
    }
}
